$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record needs to be inserted as row 238 (pushing the
# existing rows 238-310 down to 239-311). Insert a new row at position 238,
# which shifts everything below it down by one and extends the used range
# to A1:R311.
$ws.Rows.Item(238).Insert()

# Populate the newly inserted row 238 with the new data record.
$row = 238
$ws.Cells.Item($row, 1).Value2  = 8
$ws.Cells.Item($row, 2).Value2  = "Terminal La Palmera de La Serena"
$ws.Cells.Item($row, 3).Value2  = "Coquimbo"
$ws.Cells.Item($row, 4).Value2  = 45215
$ws.Cells.Item($row, 5).Value2  = 4
$ws.Cells.Item($row, 6).Value2  = 100112001
$ws.Cells.Item($row, 7).Value2  = "Berenjena"
$ws.Cells.Item($row, 8).Value2  = "Sin especificar"
$ws.Cells.Item($row, 9).Value2  = "Primera"
$ws.Cells.Item($row, 10).Value2 = 480
$ws.Cells.Item($row, 11).Value2 = 9000
$ws.Cells.Item($row, 12).Value2 = 10000
$ws.Cells.Item($row, 13).Value2 = 9500
$ws.Cells.Item($row, 14).Value2 = "`$/caja 50 unidades"
$ws.Cells.Item($row, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item($row, 16).Value2 = 190
$ws.Cells.Item($row, 17).Value2 = 50
$ws.Cells.Item($row, 18).Value2 = "Hortaliza"

# Make sure the date cell keeps the same date/time number format used by
# the rest of the date column.
$ws.Cells.Item($row, 4).NumberFormat = $ws.Cells.Item(237, 4).NumberFormat
